$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 updates
$ws.Range("G6").Value = 1.39
$ws.Range("H6").Value = 4.35
$ws.Range("I6").Value = 6.3
$ws.Range("J6").Value = 1.85
$ws.Range("K6").Value = 2.37
$ws.Range("M6").Value = 1.01
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4.05
$ws.Range("Q6").Value = 1.57
$ws.Range("R6").Value = 2.12
$ws.Range("S6").Value = 1.29
$ws.Range("T6").Value = 3.28
$ws.Range("W6").Value = 6.7
$ws.Range("X6").Value = 6.1
$ws.Range("Y6").Value = 7.1
$ws.Range("Z6").Value = 7.9
$ws.Range("AA6").Value = 9.25
$ws.Range("AB6").Value = 19
$ws.Range("AC6").Value = 13.5
$ws.Range("AD6").Value = 7.7
$ws.Range("AE6").Value = 14.5
$ws.Range("AF6").Value = 55
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 32
$ws.Range("AK6").Value = 100
$ws.Range("AM6").Value = 45
$ws.Range("AO6").Value = 6.3
$ws.Range("AP6").Value = 15.5
$ws.Range("AQ6").Value = 17.5
$ws.Range("AS6").Value = 200
$ws.Range("AU6").Value = 7.9
$ws.Range("AV6").Value = 70
$ws.Range("AW6").Value = 7.8
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 450

# Row 7 updates
$ws.Range("G7").Value = 2.57
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.62
$ws.Range("J7").Value = 3.05
$ws.Range("L7").Value = 3.25
$ws.Range("R7").Value = 1.93
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 10.5
$ws.Range("X7").Value = 15.5
$ws.Range("Y7").Value = 9.25
$ws.Range("Z7").Value = 32
$ws.Range("AA7").Value = 19.5
$ws.Range("AF7").Value = 40
$ws.Range("AH7").Value = 9.75
$ws.Range("AI7").Value = 14.5
$ws.Range("AJ7").Value = 9.5
$ws.Range("AK7").Value = 32
$ws.Range("AL7").Value = 21
$ws.Range("AM7").Value = 25
$ws.Range("AN7").Value = 4.65
$ws.Range("AO7").Value = 13.5
$ws.Range("AP7").Value = 18
$ws.Range("AQ7").Value = 55
$ws.Range("AR7").Value = 75
$ws.Range("AT7").Value = 2.8
$ws.Range("AW7").Value = 4.7
$ws.Range("AX7").Value = 14.5
